$d = $word.ActiveDocument

$replacements = @(
    @("198×7=", "461×2="),
    @("360×2=", "118×5="),
    @("201×4=", "989×3="),
    @("892×5=", "919×4="),
    @("304×5=", "951×7="),
    @("392×5=", "728×3="),
    @("104×7=", "163×3="),
    @("269×2=", "512×9="),
    @("852×2=", "483×3="),
    @("148×4=", "762×2="),
    @("507×7=", "731×8="),
    @("107×4=", "748×7="),
    @("819×8=", "427×2="),
    @("926×2=", "224×4="),
    @("953×4=", "436×6="),
    @("126×4=", "988×8="),
    @("120×2=", "609×3="),
    @("245×5=", "138×8="),
    @("856×9=", "334×7="),
    @("550×4=", "622×7="),
    @("516×8=", "334×2="),
    @("901×2=", "791×5="),
    @("672×2=", "993×8="),
    @("280×5=", "895×5="),
    @("389×9=", "856×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
